$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 0.9948644737664785
$ws.Range("C3").Value = 0.9948706811880587
$ws.Range("D3").Value = 0.9949483789308121

# Row 4 - GradientBoostingRegressor
$ws.Range("B4").Value = 0.9964251009341101
$ws.Range("C4").Value = 0.9960804471609069
$ws.Range("D4").Value = 0.9960804471609069

# Row 5 - AdaBoostRegressor
$ws.Range("B5").Value = 0.9860411958602002
$ws.Range("C5").Value = 0.9849452547299317
$ws.Range("D5").Value = 0.9836006454923457
